$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Drop the "Popis" / ${description} column (old column C) entirely and
#    rename "Nazov Hry" (old column B) to just "Nazov". Deleting column C
#    shifts everything after it one column to the left, so what was D..K
#    becomes C..J.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).Delete()
$ws.Range("B1").Value = "Názov"

# ---------------------------------------------------------------------------
# 2. Append a new trailing column (now K, since the sheet is 10 columns wide
#    after the delete) for the tutorial-video link: "Videonávod" header with
#    a ${tutorialUrl} placeholder in the data row.
# ---------------------------------------------------------------------------
$ws.Range("K1").Value = "Videonávod"
$ws.Range("K2").Value = "`${tutorialUrl}"

# ---------------------------------------------------------------------------
# 3. Re-apply explicit column widths (matches the template's refreshed
#    layout after the column removal/addition). The sheet's default column
#    width also moved (11.5703125 -> 11.640625), so column F ("Vek", which
#    inherited stale width metadata from the pre-delete column G during the
#    shift) is reset back to that new default explicitly.
# ---------------------------------------------------------------------------
$pxOffset = 0.8333333333333334
$newDefaultWidth = 11.640625
$ws.Columns.Item(3).ColumnWidth  = 18.89 - $pxOffset
$ws.Columns.Item(4).ColumnWidth  = 19.17 - $pxOffset
$ws.Columns.Item(5).ColumnWidth  = 20.42 - $pxOffset
$ws.Columns.Item(6).ColumnWidth  = $newDefaultWidth - $pxOffset
$ws.Columns.Item(7).ColumnWidth  = 15    - $pxOffset
$ws.Columns.Item(8).ColumnWidth  = 13.47 - $pxOffset
$ws.Columns.Item(9).ColumnWidth  = 14.88 - $pxOffset
$ws.Columns.Item(10).ColumnWidth = 51.96 - $pxOffset
$ws.Columns.Item(11).ColumnWidth = 50.57 - $pxOffset

# ---------------------------------------------------------------------------
# 4. Turn on the AutoFilter over the data columns (C1:J2) and register the
#    matching hidden, sheet-scoped _FilterDatabase defined name Excel uses
#    to persist it.
# ---------------------------------------------------------------------------
$ws.Range("C1:J2").AutoFilter() | Out-Null

$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$C`$1:`$J`$2")
$fdb.Visible = $false
